$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 871, shifting existing rows 871:912 down to 872:913.
$ws.Rows.Item(871).Insert()

# Populate the newly inserted row with the new data point.
# Column A holds the date as plain text ("YYYY/MM/DD"), like every other row
# in this column. Assigning the literal string straight to .Value would make
# Excel auto-detect it as a date and store a date serial instead (changing
# both the stored value and the cell format/style). Writing it as a text
# formula first, then collapsing the formula to its value via copy/paste,
# keeps it a plain text string - exactly matching the rest of the column -
# without touching the cell's style.
$ws.Cells.Item(871, 1).Formula = "=""2026/02/24"""
$ws.Cells.Item(871, 1).Copy()
$ws.Cells.Item(871, 1).PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Cells.Item(871, 2).Value = "火"
$ws.Cells.Item(871, 3).Value = 6
$ws.Cells.Item(871, 4).Value = 201
